# Apply the 30/12/2025 13:57:xx scrape update to the three worksheets:
#   - LP1912      (columns: A=Header/blank, B=Hora_Scrap, C=Hora_Llegada, D=Linea, E=Minutos, F=Parada, G=Fecha)
#   - LP1912-215  (columns: A=Header/blank, B=Fecha, C=Hora_Scrap, D=Hora_Llegada, E=Linea, F=Minutos, G=Parada)
#   - 6203-6173   (columns: A=Header/blank, B=Fecha, C=Hora_Scrap, D=Hora_Llegada, E=Linea, F=Minutos, G=Parada)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": append rows 262-273, bump header metadata
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 13:57:11"
$ws1.Range("A3").Value = "Total filas: 272"

$sheet1Rows = @(
    @(262, "13:57:01", "14:04", "17_ROMERO",            7, "LP1912", "30/12/2025"),
    @(263, "13:57:01", "14:05", "23_HERNANDEZ",         8, "LP1912", "30/12/2025"),
    @(264, "13:57:01", "14:07", "16_SANTA ANA",        10, "LP1912", "30/12/2025"),
    @(265, "13:57:01", "14:17", "16_SANTA ANA",        20, "LP1912", "30/12/2025"),
    @(266, "13:57:01", "14:21", "26_HERNANDEZ",        24, "LP1912", "30/12/2025"),
    @(267, "13:57:01", "14:45", "14_ABASTO",           48, "LP1912", "30/12/2025"),
    @(268, "13:57:01", "14:57", "16_P MOR-SANTA ANA",  60, "LP1912", "30/12/2025"),
    @(269, "13:57:01", "14:58", "215B_EL PATO",        61, "LP1912", "30/12/2025"),
    @(270, "13:57:01", "15:00", "81_EL PELIGRO",       63, "LP1912", "30/12/2025"),
    @(271, "13:57:01", "15:05", "10_OLMOS",            68, "LP1912", "30/12/2025"),
    @(272, "13:57:01", "15:20", "15_ABASTO",           83, "LP1912", "30/12/2025"),
    @(273, "13:57:01", "15:23", "26_HERNANDEZ",        86, "LP1912", "30/12/2025")
)

foreach ($r in $sheet1Rows) {
    $rowNum = $r[0]
    $ws1.Cells.Item($rowNum, 2).Value = $r[1]
    $ws1.Cells.Item($rowNum, 3).Value = $r[2]
    $ws1.Cells.Item($rowNum, 4).Value = $r[3]
    $ws1.Cells.Item($rowNum, 5).Value = $r[4]
    $ws1.Cells.Item($rowNum, 6).Value = $r[5]
    $ws1.Cells.Item($rowNum, 7).Value = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": append row 21, bump header metadata
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 13:57:11"
$ws2.Range("A3").Value = "Total filas: 20"

$ws2.Cells.Item(21, 2).Value = "30/12/2025"
$ws2.Cells.Item(21, 3).Value = "13:57:01"
$ws2.Cells.Item(21, 4).Value = "14:58"
$ws2.Cells.Item(21, 5).Value = "215B_EL PATO"
$ws2.Cells.Item(21, 6).Value = 61
$ws2.Cells.Item(21, 7).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": append rows 39-41, bump header metadata
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 13:57:11"
$ws3.Range("A3").Value = "Total filas: 40"

$sheet3Rows = @(
    @(39, "30/12/2025", "13:57:11", "14:09", "215A_LA PLATA", 12, "L6173"),
    @(40, "30/12/2025", "13:57:06", "14:53", "215D_LA PLATA", 56, "L6203"),
    @(41, "30/12/2025", "13:57:11", "15:34", "215A_LA PLATA", 97, "L6173")
)

foreach ($r in $sheet3Rows) {
    $rowNum = $r[0]
    $ws3.Cells.Item($rowNum, 2).Value = $r[1]
    $ws3.Cells.Item($rowNum, 3).Value = $r[2]
    $ws3.Cells.Item($rowNum, 4).Value = $r[3]
    $ws3.Cells.Item($rowNum, 5).Value = $r[4]
    $ws3.Cells.Item($rowNum, 6).Value = $r[5]
    $ws3.Cells.Item($rowNum, 7).Value = $r[6]
}
